$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('3153 Sacred Heart Community St Kilda Tier 1B Outbreak', 15, 'new'),
    @('3175 The Bays Aged Care Facility Hastings', 16, 'old'),
    @('3175 The Bays Aged Care Facility Hastings Outbreak', 15, 'new'),
    @('3563 Embracia Aged Care Reservoir', 11, 'old'),
    @('3600 Belvedere Aged Care Noble Park', 22, 'old'),
    @('3600 Belvedere Aged Care Noble Park Outbreak', 27, 'new'),
    @('3612 BlueCross Glengowrie', 35, 'old'),
    @('3612 BlueCross Glengowrie Outbreak', 35, 'new'),
    @('3684 Homestyle Aged Care Langford Grange Cranbourne East', 25, 'old'),
    @('3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak', 31, 'new'),
    @('3980 Arcare Keysborough Aged Care Keysborough', 13, 'old'),
    @('3980 Arcare Keysborough Aged Care Keysborough Outbreak', 19, 'new'),
    @('4075 Ferndale Gardens Aged Care Services Bayswater North', 17, 'old'),
    @('4075 Ferndale Gardens Aged Care Services Bayswater North Outbreak', 17, 'new'),
    @('4518 Regis Aged Care Fawkner', 14, 'old'),
    @('4518 Regis Aged Care Fawkner Outbreak', 13, 'new'),
    @('AW Window Transport Group Depot North Geelong', 10, 'old'),
    @('AW Window Transport Group Depot North Geelong Outbreak', 10, 'new'),
    @('Allied Pinnacle Factory Altona North', 12, 'old'),
    @('Allied Pinnacle Factory Altona North Outbreak', 15, 'new'),
    @('Bespoke Childcare Dingley Village Outbreak', 10, 'new'),
    @('Bread Solutions Braeside', 15, 'old'),
    @('Bread Solutions Braeside Outbreak', 15, 'new'),
    @('CS Square Caroline Springs', 13, 'old'),
    @('CS Square Caroline Springs Outbreak', 11, 'new'),
    @('Cedar Meats Australia Brooklyn', 10, 'old'),
    @('Child''s Play Early Learning Centre Tarneit', 11, 'old'),
    @('Child''s Play Early Learning Centre Tarneit Outbreak', 11, 'new'),
    @('Comfort Sleep Bedding Co Thomastown', 10, 'old'),
    @('Community Kids Pascoe Vale Early Education Centre Pascoe Vale', 22, 'old'),
    @('Community Kids Pascoe Vale Early Education Centre Pascoe Vale Outbreak', 22, 'new'),
    @('Essential Caravans Somerton Outbreak', 10, 'new'),
    @('Guardian Childcare Caulfield', 20, 'old'),
    @('Guardian Childcare Caulfield Outbreak', 21, 'new'),
    @('Hello Fresh Warehouse Ravenhall', 17, 'old'),
    @('Hello Fresh Warehouse Ravenhall Outbreak', 17, 'new'),
    @('Inghams Enterprises Somerville', 24, 'old'),
    @('Inghams Enterprises Somerville Outbreak', 24, 'new'),
    @('Lantmannen Unibake Australia Mordialloc', 22, 'old'),
    @('Lantmannen Unibake Australia Mordialloc Outbreak', 21, 'new'),
    @('Launch Housing City Edge Crisis Accommodation South Melbourne', 11, 'old'),
    @('Launch Housing City Edge Crisis Accommodation South Melbourne Outbreak', 11, 'new'),
    @('MacKillop Family Services Residential Facility Glenroy Outbreak', 10, 'new'),
    @('Melbourne Assessment Prison West Melbourne Outbreak', 16, 'new'),
    @('Melbourne Youth Justice Centre Parkville Outbreak', 10, 'new'),
    @('Mildura Police Station Mildura Outbreak', 10, 'new'),
    @('Monash Health Kingston Centre South 5', 10, 'old'),
    @('Monash Health Kingston Centre South 5 Outbreak', 10, 'new'),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 41, 'old'),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B Outbreak', 43, 'new'),
    @('Northern Health The Northern Hospital Epping', 17, 'old'),
    @('Northern Health The Northern Hospital Epping Outbreak', 19, 'new'),
    @('Oceania Meat Processors Laverton North', 16, 'old'),
    @('Oceania Meat Processors Laverton North Outbreak', 17, 'new'),
    @('Public Housing 140 Brunswick Street Fitzroy Outbreak', 11, 'new'),
    @('Shawlands Caravan Park Dandenong South Outbreak', 17, 'new'),
    @('St Vincents Hospital Emergency Department Melbourne', 36, 'old'),
    @('St Vincents Hospital Emergency Department Melbourne Outbreak', 37, 'new'),
    @('Target Distribution Centre Truganina', 17, 'old'),
    @('Target Distribution Centre Truganina Outbreak', 11, 'new'),
    @('The Robin Hood Inn Drouin West', 46, 'old'),
    @('The Robin Hood Inn Drouin West Outbreak', 48, 'new'),
    @('The Royal Melbourne Hospital Parkville Emergency Department', 10, 'old'),
    @('The Toolshed Bar Private Event Noojee', 16, 'old'),
    @('The Toolshed Bar Private Event Noojee Outbreak', 16, 'new'),
    @('Turosi Breakwater', 10, 'old'),
    @('Turosi Breakwater Outbreak', 10, 'new'),
    @('Visy Recycling Springvale', 21, 'old'),
    @('Visy Recycling Springvale Outbreak', 14, 'new'),
    @('Werribee Mercy Hospital Emergency Department', 23, 'old'),
    @('Werribee Mercy Hospital Emergency Department Outbreak', 26, 'new'),
    @('Western Health Footscray Hospital Emergency Department', 10, 'old'),
    @('Western Health Footscray Hospital Emergency Department Outbreak', 11, 'new'),
    @('Western Health Sunshine Hospital Emergency Department', 22, 'old'),
    @('Western Health Sunshine Hospital Emergency Department Outbreak', 24, 'new'),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}